$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Implementação avançada: funcionalidades de  frequência, notificações, upload de fotos, treino, avaliação física. Refinamento do protótipo. Conexão da tela de criação de instrutor com a API"

$ws.Range("A12").Value = "Continuação da implementação avaçada, incluindo a criação das telas do usuário do tipo Instrutor. Revisão e ajustes de interface e de código conforme feedbacks, tanto dos stackholders como dos professores.`nPreparação de uma nova versão do documento de requisitos e refinamento dos protótipos de tela."

$ws.Range("A13").Value = "Continuação da implementação avançada. Continuação dos testes do sistema.`nFinalização do documentação de requisitos."

$ws.Range("A14").Value = "Contiuação da implemetnação: ajustes finais nas funcionalidades e documentação.`nPreparação do sistema para avaliação com potenciais usuários."
